# Update the pension-recipient figures for 2015-2021 (row 4, columns E:K)
# on the single worksheet, per the commit "files updated and bug fixed".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E4").Value = 7293
$ws.Range("F4").Value = 7384
$ws.Range("G4").Value = 7532
$ws.Range("H4").Value = 7652
$ws.Range("I4").Value = 7801
$ws.Range("J4").Value = 8086
$ws.Range("K4").Value = 8158
